$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''303.71'
$ws.Range("E2").Value = '''5.26%'
$ws.Range("D3").Value = '''34.85'
$ws.Range("E3").Value = '''12.32%'
$ws.Range("D4").Value = '''5.131'
$ws.Range("E4").Value = '''4.30%'
$ws.Range("D5").Value = '''0.07760'
$ws.Range("E5").Value = '''5.16%'
$ws.Range("D6").Value = '''2.365'
$ws.Range("E6").Value = '''6.76%'
$ws.Range("D7").Value = '''8.028'
$ws.Range("E7").Value = '''4.25%'
$ws.Range("D8").Value = '''3.933'
$ws.Range("E8").Value = '''5.41%'
$ws.Range("D9").Value = '''0.9298'
$ws.Range("E9").Value = '''2.21%'
$ws.Range("E10").Value = '''16.33%'
$ws.Range("D11").Value = '''0.1796'
$ws.Range("E11").Value = '''6.09%'
$ws.Range("E12").Value = '''3.88%'
$ws.Range("D13").Value = '''0.03311'
$ws.Range("E13").Value = '''6.47%'
$ws.Range("D14").Value = '''0.09896'
$ws.Range("E14").Value = '''-0.47%'
$ws.Range("D15").Value = '''0.001498'
$ws.Range("E15").Value = '''-0.10%'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04565'
$ws.Range("E16").Value = '''0.38%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.005777'
$ws.Range("E17").Value = '''-0.60%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.466'
$ws.Range("E18").Value = '''-0.73%'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '''2.133'
$ws.Range("E19").Value = '''2.02%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3367'
$ws.Range("E20").Value = '''1.18%'
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '''0.1310'
$ws.Range("E21").Value = '''1.12%'
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D22").Value = '''4.308'
$ws.Range("E22").Value = '''12.68%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.2387'
$ws.Range("E23").Value = '''12.49%'
$ws.Range("E24").Value = '''0.50%'
$ws.Range("D25").Value = '''0.004461'
$ws.Range("E25").Value = '''7.65%'
$ws.Range("D26").Value = '''0.0001299'
$ws.Range("E26").Value = '''-0.20%'
$ws.Range("E27").Value = '''-0.10%'
$ws.Range("D39").Value = '''0.01782'
$ws.Range("E39").Value = '''12.21%'
$ws.Range("D40").Value = '''0.04760'
$ws.Range("E40").Value = '''6.50%'
$ws.Range("D41").Value = '''0.007755'
$ws.Range("E41").Value = '''5.63%'
$ws.Range("D42").Value = '''0.1413'
$ws.Range("E42").Value = '''6.67%'
$ws.Range("D43").Value = '''0.007078'
$ws.Range("E43").Value = '''-25.69%'
$ws.Range("D44").Value = '''0.002148'
$ws.Range("E44").Value = '''-3.35%'
$ws.Range("D45").Value = '''0.009194'
$ws.Range("E45").Value = '''9.26%'
$ws.Range("E46").Value = '''0.10%'
$ws.Range("E47").Value = '''-0.05%'
$ws.Range("D48").Value = '''2.741'
$ws.Range("E48").Value = '''21.47%'
$ws.Range("E49").Value = '''-0.05%'
$ws.Range("E50").Value = '''-0.05%'
$ws.Range("E51").Value = '''-0.05%'
